# Add data for a few other cultivars in CE experiments.
# - New column H ("Wheat.SowingData.Cultivar") is populated for all existing
#   data rows (238-261) with the cultivar each simulation used.
# - 16 new data rows (262-277) are appended for four more cultivars:
#   Claire, Mccubbin, Rongotea and Wakanui (LV/LN/SV/SN durations each).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New header cell H1 : "Wheat.SowingData.Cultivar"
# Styled like the other headers (bold, bordered) but left-aligned with an
# explicit black font colour, matching the rest of row 1's look.
# ---------------------------------------------------------------------------
$h1 = $ws.Range("H1")
$h1.Value = "Wheat.SowingData.Cultivar"
$h1.Font.Bold = $true
$h1.Font.Color = 0
$h1.HorizontalAlignment = -4131   # xlLeft
$h1.VerticalAlignment = -4160     # xlTop
$h1.Borders.LineStyle = 1
$h1.Borders.Weight = 2

# ---------------------------------------------------------------------------
# Populate the new Cultivar column for the existing rows (238-261)
# ---------------------------------------------------------------------------
$cultivars238to261 = @(
    "Amarok","Amarok",
    "CRWT153","CRWT153",
    "Otane","Otane",
    "Saracen","Saracen",
    "BattenWinter","BattenWinter",
    "BattenSpring","BattenSpring",
    "Amarok","Amarok",
    "CRWT153","CRWT153",
    "Otane","Otane",
    "Saracen","Saracen",
    "BattenWinter","BattenWinter",
    "BattenSpring","BattenSpring"
)

$startRow = 238
for ($i = 0; $i -lt $cultivars238to261.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 8).Value = $cultivars238to261[$i]
}

# ---------------------------------------------------------------------------
# New rows 262-277 : four extra cultivars, each with LV / LN / SV / SN
# HarvestRipe duration records.
# Columns: A = SimulationName, C = CurrentStageName, D = duration, H = Cultivar
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 262; A = "LincolnCETreatLVCvClaireDurat12";   D = 9;   H = "Claire" },
    @{ Row = 263; A = "LincolnCETreatLNCvClaireDurat12";   D = 16;  H = "Claire" },
    @{ Row = 264; A = "LincolnCETreatSVCvClaireDurat12";   D = 12;  H = "Claire" },
    @{ Row = 265; A = "LincolnCETreatSNCvClaireDurat12";   D = 13;  H = "Claire" },
    @{ Row = 266; A = "LincolnCETreatLVCvMccubbinDurat12"; D = 8;   H = "Mccubbin" },
    @{ Row = 267; A = "LincolnCETreatLNCvMccubbinDurat12"; D = 8;   H = "Mccubbin" },
    @{ Row = 268; A = "LincolnCETreatSVCvMccubbinDurat12"; D = 16;  H = "Mccubbin" },
    @{ Row = 269; A = "LincolnCETreatSNCvMccubbinDurat12"; D = 13;  H = "Mccubbin" },
    @{ Row = 270; A = "LincolnCETreatLVCvRongoteaDurat12"; D = 8;   H = "Rongotea" },
    @{ Row = 271; A = "LincolnCETreatLNCvRongoteaDurat12"; D = 7;   H = "Rongotea" },
    @{ Row = 272; A = "LincolnCETreatSVCvRongoteaDurat12"; D = 9;   H = "Rongotea" },
    @{ Row = 273; A = "LincolnCETreatSNCvRongoteaDurat12"; D = 11;  H = "Rongotea" },
    @{ Row = 274; A = "LincolnCETreatLVCvWakanuiDurat12";  D = 9.5; H = "Wakanui" },
    @{ Row = 275; A = "LincolnCETreatLNCvWakanuiDurat12";  D = 15;  H = "Wakanui" },
    @{ Row = 276; A = "LincolnCETreatSVCvWakanuiDurat12";  D = 11;  H = "Wakanui" },
    @{ Row = 277; A = "LincolnCETreatSNCvWakanuiDurat12";  D = 17;  H = "Wakanui" }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 3).Value = "HarvestRipe"
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 8).Value = $r.H
}

# Style column A for the new rows the same way the rest of column A is
# styled (bold font, top-aligned & centred).
# Rows 262-265 (Claire) keep the full box border used throughout column A.
$claireA = $ws.Range("A262:A265")
$claireA.Font.Bold = $true
$claireA.HorizontalAlignment = -4108   # xlCenter
$claireA.VerticalAlignment = -4160     # xlTop
$claireA.Borders.LineStyle = 1
$claireA.Borders.Weight = 2

# Rows 266-277 (Mccubbin, Rongotea, Wakanui) use a narrower left/right-only
# border variant.
$restA = $ws.Range("A266:A277")
$restA.Font.Bold = $true
$restA.HorizontalAlignment = -4108   # xlCenter
$restA.VerticalAlignment = -4160     # xlTop
$restA.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$restA.Borders.Item(7).Weight = 2
$restA.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$restA.Borders.Item(10).Weight = 2

# ---------------------------------------------------------------------------
# Match the author's final selection / scroll position
# ---------------------------------------------------------------------------
$null = $ws.Range("H1").Select()
